$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data: keyword logic-based dialog act classifier, 2 hours in F
$ws.Range("A3").Value = "keyword logic-based dialog act classifier"
$ws.Range("F3").Value = 2

# Adjust column A width to fit the new, longer text
$ws.Columns.Item(1).ColumnWidth = 33.5

# Update selection to reflect where the user ended up (G15)
$ws.Range("G15").Select()
